$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.936.31"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.301.37"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "300.38"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "97.31"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").Value = "35.79"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D12").Value = "17.89"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "2.656.21"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "2.306.13"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "42.880.26"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").Value = "67.87"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "240.80"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "25.46"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "165.72"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").Value = "33.03"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D36").Value = "17.10"
$ws.Range("E36").Value = "  -6.82%  "
$ws.Range("D37").Value = "2.37"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").Value = "0.0686"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D43").Value = "2.017.63"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").Value = "10.17"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "2.14"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").Value = "17.39"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "2.94"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").Value = "53.51"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "2.522.10"
$ws.Range("E51").Value = "  -0.68%  "
